$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 45877.4169196875

$ws.Range("A12").Value = 45877.45849543768
$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 32
$ws.Range("D12").Value = 15.44
$ws.Range("E12").Value = 89.29000000000001
$ws.Range("F12").Value = 302.34
$ws.Range("G12").Value = 9.31
$ws.Range("H12").Value = "ESE"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "11:00:14"

$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
